# Update Castelvetro.xlsx report: append new daily rows (through 2021-12-08 / 8/12)
# as described by the commit "aggiornamento fino a 8/12".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 386-464 (serial date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44460,1,15,132.8727079457879),
    @(44461,0,15,132.8727079457879),
    @(44462,1,13,115.1563468863496),
    @(44463,1,13,115.1563468863496),
    @(44464,2,9,79.72362476747276),
    @(44465,0,5,44.29090264859597),
    @(44466,0,5,44.29090264859597),
    @(44467,0,4,35.43272211887678),
    @(44468,0,4,35.43272211887678),
    @(44469,0,3,26.57454158915759),
    @(44470,2,4,35.43272211887678),
    @(44471,0,2,17.71636105943839),
    @(44472,0,2,17.71636105943839),
    @(44473,0,2,17.71636105943839),
    @(44474,1,3,26.57454158915759),
    @(44475,1,4,35.43272211887678),
    @(44476,5,9,79.72362476747276),
    @(44477,0,7,62.00726370803437),
    @(44478,0,7,62.00726370803437),
    @(44479,0,7,62.00726370803437),
    @(44480,0,7,62.00726370803437),
    @(44481,0,6,53.14908317831517),
    @(44482,0,5,44.29090264859597),
    @(44483,0,0,0),
    @(44484,0,0,0),
    @(44485,0,0,0),
    @(44486,0,0,0),
    @(44487,0,0,0),
    @(44488,0,0,0),
    @(44489,0,0,0),
    @(44490,0,0,0),
    @(44491,0,0,0),
    @(44492,0,0,0),
    @(44493,0,0,0),
    @(44494,0,0,0),
    @(44495,0,0,0),
    @(44496,0,0,0),
    @(44497,1,1,8.858180529719196),
    @(44498,1,2,17.71636105943839),
    @(44499,2,4,35.43272211887678),
    @(44500,0,4,35.43272211887678),
    @(44501,0,4,35.43272211887678),
    @(44502,0,4,35.43272211887678),
    @(44503,0,4,35.43272211887678),
    @(44504,0,3,26.57454158915759),
    @(44505,0,2,17.71636105943839),
    @(44506,2,2,17.71636105943839),
    @(44507,1,3,26.57454158915759),
    @(44508,3,6,53.14908317831517),
    @(44509,1,7,62.00726370803437),
    @(44510,0,7,62.00726370803437),
    @(44511,3,10,88.58180529719195),
    @(44512,0,10,88.58180529719195),
    @(44513,1,9,79.72362476747276),
    @(44514,3,11,97.43998582691115),
    @(44515,7,15,132.8727079457879),
    @(44516,7,21,186.0217911241031),
    @(44517,0,21,186.0217911241031),
    @(44518,0,18,159.4472495349455),
    @(44519,3,21,186.0217911241031),
    @(44520,1,21,186.0217911241031),
    @(44521,2,20,177.1636105943839),
    @(44522,2,15,132.8727079457879),
    @(44523,6,14,124.0145274160687),
    @(44524,5,19,168.3054300646647),
    @(44525,4,23,203.7381521835415),
    @(44526,7,27,239.1708743024183),
    @(44527,4,30,265.7454158915759),
    @(44528,7,35,310.0363185401719),
    @(44529,9,42,372.0435822482062),
    @(44530,3,39,345.4690406590486),
    @(44531,1,35,310.0363185401719),
    @(44532,22,53,469.4835680751174),
    @(44533,5,51,451.7672070156789),
    @(44534,8,55,487.1999291345558),
    @(44535,3,51,451.7672070156789),
    @(44536,6,48,425.1926654265214),
    @(44537,7,52,460.6253875453982),
    @(44538,2,53,469.4835680751174)
)


$startRow = 386
$rowCount = $data.Count
$colCount = 4

# Build a 2D array for a single bulk write (serial date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$targetRange.Value2 = $arr

# Copy the formatting (date number format, bold font, thin border, center/top alignment)
# of the last pre-existing date cell (A385) down over the newly added date cells,
# matching the same per-row style already used for column A.
$srcFormat = $ws.Range("A385")
$dstFormat = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)
$excel.CutCopyMode = 0
